$wb = $excel.ActiveWorkbook

# --- ODI Batting: B3 was an empty inline-string placeholder cell, drop it ---
$wsBatting = $wb.Worksheets.Item("ODI Batting")
$wsBatting.Range("B3").Value = $null
$wsBatting.Range("E3").Value = " "

# --- Add new sheet "ODI Batting Extra" right after "ODI Batting" ---
$wsExtra = $wb.Worksheets.Add($null, $wsBatting)
$wsExtra.Name = "ODI Batting Extra"

# Header row (bold, centered, top-aligned, thin border - matches the other sheets' header style)
$headers = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $wsExtra.Cells.Item(1, $i + 1)
    $cell.Value = $headers[$i]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
    $cell.Borders.Weight = 2
}

# Row 2
$wsExtra.Range("A2").Value = "'4401"
$wsExtra.Range("F2").Value = "NO"

# Row 3
$wsExtra.Range("A3").Value = "'4405"
$wsExtra.Range("F3").Value = "NO"

# Row 4
$wsExtra.Range("A4").Value = "'4408"
$wsExtra.Range("F4").Value = "NO"

# Row 5
$wsExtra.Range("A5").Value = "'4426"
$wsExtra.Range("F5").Value = "NO"

# Row 6
$wsExtra.Range("A6").Value = "'4427"
$wsExtra.Range("B6").Value = 4
$wsExtra.Range("C6").Value = "'3"
$wsExtra.Range("D6").Value = "'0"
$wsExtra.Range("E6").Value = "'6.94%"
$wsExtra.Range("F6").Value = "NO"

# Row 7
$wsExtra.Range("A7").Value = "'4428"
$wsExtra.Range("B7").Value = 5
$wsExtra.Range("C7").Value = "'6"
$wsExtra.Range("D7").Value = "'1"
$wsExtra.Range("E7").Value = "'17.68%"
$wsExtra.Range("F7").Value = "NO"
